$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$colA = $ws.Range("A:A")

# 1. Remove the THAIS row (account 005395948, balance 50000) entirely.
$thais = $colA.Find("005395948")
$ws.Rows($thais.Row()).Delete()

# 2. Remove the old JULIANA row (account 004813088, balance 38195.84) - it will
#    be re-inserted further up with its new balance.
$julianaOld = $colA.Find("004813088")
$ws.Rows($julianaOld.Row()).Delete()

# 3. Insert the JULIANA row (new balance 54565.49) right before RENATO
#    (account 000330949), i.e. directly after EDUARDO.
# Leading "'" forces text type so the leading zeros of the account number survive.
$renato = $colA.Find("000330949")
$ws.Rows($renato.Row()).Insert()
$ws.Cells.Item($renato.Row(), 1).Value = "'004813088"
$ws.Cells.Item($renato.Row(), 2).Value = "JULIANA"
$ws.Cells.Item($renato.Row(), 3).Value = 54565.49

# 4. Insert the BLUEMETRIX row right before CARLOS (account 004488571).
$carlos = $colA.Find("004488571")
$ws.Rows($carlos.Row()).Insert()
$ws.Cells.Item($carlos.Row(), 1).Value = "'001761119"
$ws.Cells.Item($carlos.Row(), 2).Value = "BLUEMETRIX"
$ws.Cells.Item($carlos.Row(), 3).Value = 600.26

# 5. Insert the EDMUR row right after THEOMAR (account 004382374, balance 129.48).
$theomar = $colA.Find("004382374")
$edmurRow = $theomar.Row() + 1
$ws.Rows($edmurRow).Insert()
$ws.Cells.Item($edmurRow, 1).Value = "'005009026"
$ws.Cells.Item($edmurRow, 2).Value = "EDMUR"
$ws.Cells.Item($edmurRow, 3).Value = 109.45
